$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.982.29'
$ws.Range('E2').Value = '  -3.99%  '
$ws.Range('D3').Value = '1.642.75'
$ws.Range('E3').Value = '  -3.33%  '
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '214.90'
$ws.Range('E5').Value = '  -3.66%  '
$ws.Range('D6').Value = '0.5077'
$ws.Range('E6').Value = '  -3.08%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = '0.2576'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('D9').Value = '0.06399'
$ws.Range('E9').Value = '  -2.92%  '
$ws.Range('D10').Value = '19.54'
$ws.Range('E10').Value = '  -4.38%  '
$ws.Range('D11').Value = '0.07729'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '1.646.36'
$ws.Range('E12').Value = '  -3.33%  '
$ws.Range('D13').Value = '4.252'
$ws.Range('E13').Value = '  -3.49%  '
$ws.Range('D14').Value = '1.871.61'
$ws.Range('E14').Value = '  -3.19%  '
$ws.Range('D15').Value = '0.5436'
$ws.Range('E15').Value = '  -4.74%  '
$ws.Range('D16').Value = '0.0₅7945'
$ws.Range('E16').Value = '  -1.89%  '
$ws.Range('D17').Value = '63.60'
$ws.Range('E17').Value = '  -4.39%  '
$ws.Range('D18').Value = '26.009.21'
$ws.Range('E18').Value = '  -4.09%  '
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = '204.99'
$ws.Range('E20').Value = '  -5.65%  '
$ws.Range('D21').Value = '4.352'
$ws.Range('E21').Value = '  -4.78%  '
$ws.Range('D22').Value = '9.992'
$ws.Range('E22').Value = '  -2.66%  '
$ws.Range('D23').Value = '5.977'
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = '1.923'
$ws.Range('E25').Value = '  +11.31%  '
$ws.Range('D26').Value = '143.19'
$ws.Range('E26').Value = '  -1.09%  '
$ws.Range('D27').Value = '0.1161'
$ws.Range('E27').Value = '  -2.35%  '
$ws.Range('D28').Value = '6.870'
$ws.Range('E28').Value = '  -3.56%  '
$ws.Range('D29').Value = '15.74'
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').Value = '0.05037'
$ws.Range('E30').Value = '  -4.68%  '
$ws.Range('D31').Value = '1.237'
$ws.Range('E31').Value = '  -3.85%  '
$ws.Range('D32').Value = '3.303'
$ws.Range('E32').Value = '  -3.31%  '
$ws.Range('D33').Value = '3.211'
$ws.Range('E33').Value = '  -2.75%  '
$ws.Range('D34').Value = '1.536'
$ws.Range('E34').Value = '  -5.07%  '
$ws.Range('E35').Value = '  -2.03%  '
$ws.Range('D36').Value = '0.9104'
$ws.Range('E36').Value = '  -3.43%  '
$ws.Range('D37').Value = '2.645'
$ws.Range('E37').Value = '  -5.86%  '
$ws.Range('D38').Value = '0.5673'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('D39').Value = '1.139.61'
$ws.Range('E39').Value = '  -3.37%  '
$ws.Range('D40').Value = '0.01566'
$ws.Range('E40').Value = '  -4.07%  '
$ws.Range('D41').Value = '2.560'
$ws.Range('E41').Value = '  -0.58%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Value = '5.634'
$ws.Range('E43').Value = '  -1.52%  '
$ws.Range('D44').Value = '0.8206'
$ws.Range('E44').Value = '  -1.96%  '
$ws.Range('D45').Value = '99.61'
$ws.Range('E45').Value = '  -0.88%  '
$ws.Range('D46').Value = '1.785.66'
$ws.Range('E46').Value = '  -3.06%  '
$ws.Range('D47').Value = '0.0₈114'
$ws.Range('E47').Value = '  +2.89%  '
$ws.Range('D48').Value = '0.4533'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').Value = '1.008'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').Value = '54.92'
$ws.Range('E50').Value = '  -3.39%  '
$ws.Range('D51').Value = '7.789'
$ws.Range('E51').Value = '  -3.33%  '
